$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New data rows to insert for Data sheet (rows 3..18), replacing the old
# rows 3..9 (1Y..30Y OIS) with the new Futures + OIS combination.
$data = @(
    @("3M",  "SQZ25",   "FUTURE", 96.14),
    @("5M",  "SQF26",   "FUTURE", 96.235),
    @("5M",  "SQG26",   "FUTURE", 96.315),
    @("6M",  "SQH26",   "FUTURE", 96.375),
    @("9M",  "SQM26",   "FUTURE", 96.61499999999999),
    @("12M", "SQU26",   "FUTURE", 96.78),
    @("15M", "SQZ26",   "FUTURE", 96.86499999999999),
    @("0M",  "SQU25",   "FUTURE", 95.86),
    @("1M",  "SQV25",   "FUTURE", 95.95999999999999),
    @("2M",  "SQX25",   "FUTURE", 96.065),
    @("2Y",  "SOFROIS", "OIS",    0.03537),
    @("3Y",  "SOFROIS", "OIS",    0.034488),
    @("5Y",  "SOFROIS", "OIS",    0.034999),
    @("10Y", "SOFROIS", "OIS",    0.038125),
    @("20Y", "SOFROIS", "OIS",    0.04137),
    @("30Y", "SOFROIS", "OIS",    0.04107)
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
